$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filters")

# Update the product image URLs (imageUrl column, B) for three filter rows
# to the newly-hosted mlstatic.com images used for the deploy.
$ws.Range("B11").Value = "https://http2.mlstatic.com/D_NQ_NP_2X_881694-MLA70061350186_062023-F.webp"
$ws.Range("B12").Value = "https://http2.mlstatic.com/D_NQ_NP_2X_815627-MLA72720301527_112023-F.webp"
$ws.Range("B13").Value = "https://http2.mlstatic.com/D_NQ_NP_2X_680318-MLA75641583430_042024-F.webp"

# Leave the current selection where the author ended up after editing.
$ws.Activate() | Out-Null
$ws.Range("B28").Select() | Out-Null
